$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New weather observation rows 366-388 appended below the existing data (row 365 was the prior last row).
# Each inner array is (row, col, value) for a single cell; string entries are shared-string text values
# (column Z holds weather-condition text: 'Rain' / 'Rain-Thunderstorm').
$cellData = @(
    ,(366,6,18)
    ,(366,7,11)
    ,(366,8,4)
    ,(366,9,-1)
    ,(366,10,-8)
    ,(366,11,-11)
    ,(366,12,39)
    ,(366,13,24)
    ,(366,14,6)
    ,(366,15,1017)
    ,(366,16,1014)
    ,(366,17,1011)
    ,(366,18,14)
    ,(366,19,10)
    ,(366,20,10)
    ,(366,21,37)
    ,(366,22,21)
    ,(366,24,0)
    ,(366,25,6)
    ,(366,27,276)
    ,(367,6,19)
    ,(367,7,16)
    ,(367,8,13)
    ,(367,9,3)
    ,(367,10,-1)
    ,(367,11,-2)
    ,(367,12,44)
    ,(367,13,31)
    ,(367,14,15)
    ,(367,15,1016)
    ,(367,16,1013)
    ,(367,17,1007)
    ,(367,18,10)
    ,(367,19,10)
    ,(367,20,10)
    ,(367,21,29)
    ,(367,22,10)
    ,(367,24,0)
    ,(367,25,5)
    ,(367,26,"Rain-Thunderstorm")
    ,(367,27,67)
    ,(368,6,16)
    ,(368,7,12)
    ,(368,8,8)
    ,(368,9,6)
    ,(368,10,4)
    ,(368,11,1)
    ,(368,12,82)
    ,(368,13,63)
    ,(368,14,34)
    ,(368,15,1013)
    ,(368,16,1011)
    ,(368,17,1004)
    ,(368,18,10)
    ,(368,19,9)
    ,(368,20,3)
    ,(368,21,55)
    ,(368,22,16)
    ,(368,24,8.89)
    ,(368,25,5)
    ,(368,26,"Rain-Thunderstorm")
    ,(368,27,169)
    ,(369,6,15)
    ,(369,7,9)
    ,(369,8,3)
    ,(369,9,4)
    ,(369,10,-2)
    ,(369,11,-11)
    ,(369,12,81)
    ,(369,13,46)
    ,(369,14,9)
    ,(369,15,1015)
    ,(369,16,1012)
    ,(369,17,1009)
    ,(369,18,14)
    ,(369,19,10)
    ,(369,20,8)
    ,(369,21,48)
    ,(369,22,26)
    ,(369,24,0)
    ,(369,25,1)
    ,(369,27,281)
    ,(370,6,19)
    ,(370,7,12)
    ,(370,8,6)
    ,(370,9,-2)
    ,(370,10,-5)
    ,(370,11,-8)
    ,(370,12,53)
    ,(370,13,30)
    ,(370,14,8)
    ,(370,15,1021)
    ,(370,16,1019)
    ,(370,17,1015)
    ,(370,18,11)
    ,(370,19,10)
    ,(370,20,10)
    ,(370,21,23)
    ,(370,22,10)
    ,(370,24,0)
    ,(370,25,6)
    ,(370,27,265)
    ,(371,6,21)
    ,(371,7,16)
    ,(371,8,9)
    ,(371,9,0)
    ,(371,10,-4)
    ,(371,11,-10)
    ,(371,12,44)
    ,(371,13,25)
    ,(371,14,6)
    ,(371,15,1022)
    ,(371,16,1019)
    ,(371,17,1013)
    ,(371,18,10)
    ,(371,19,10)
    ,(371,20,10)
    ,(371,21,29)
    ,(371,22,8)
    ,(371,24,0)
    ,(371,25,5)
    ,(371,27,263)
    ,(372,6,18)
    ,(372,7,14)
    ,(372,8,11)
    ,(372,9,8)
    ,(372,10,3)
    ,(372,11,-4)
    ,(372,12,77)
    ,(372,13,49)
    ,(372,14,11)
    ,(372,15,1021)
    ,(372,16,1018)
    ,(372,17,1013)
    ,(372,18,10)
    ,(372,19,8)
    ,(372,20,6)
    ,(372,21,34)
    ,(372,22,10)
    ,(372,24,2.0299999999999998)
    ,(372,25,7)
    ,(372,26,"Rain")
    ,(372,27,12)
    ,(373,6,14)
    ,(373,7,12)
    ,(373,8,10)
    ,(373,9,9)
    ,(373,10,8)
    ,(373,11,7)
    ,(373,12,88)
    ,(373,13,73)
    ,(373,14,52)
    ,(373,15,1018)
    ,(373,16,1013)
    ,(373,17,1007)
    ,(373,18,11)
    ,(373,19,8)
    ,(373,20,3)
    ,(373,21,19)
    ,(373,22,10)
    ,(373,24,3.05)
    ,(373,25,8)
    ,(373,26,"Rain")
    ,(373,27,62)
    ,(374,6,21)
    ,(374,7,16)
    ,(374,8,10)
    ,(374,9,8)
    ,(374,10,4)
    ,(374,11,-2)
    ,(374,12,82)
    ,(374,13,48)
    ,(374,14,15)
    ,(374,15,1012)
    ,(374,16,1009)
    ,(374,17,1003)
    ,(374,18,10)
    ,(374,19,10)
    ,(374,20,8)
    ,(374,21,34)
    ,(374,22,14)
    ,(374,24,0.25)
    ,(374,25,4)
    ,(374,26,"Rain")
    ,(374,27,286)
    ,(375,6,19)
    ,(375,7,15)
    ,(375,8,11)
    ,(375,9,7)
    ,(375,10,-2)
    ,(375,11,-11)
    ,(375,12,76)
    ,(375,13,33)
    ,(375,14,6)
    ,(375,15,1010)
    ,(375,16,1008)
    ,(375,17,1001)
    ,(375,18,14)
    ,(375,19,10)
    ,(375,20,7)
    ,(375,21,37)
    ,(375,22,19)
    ,(375,24,0.25)
    ,(375,25,4)
    ,(375,26,"Rain")
    ,(375,27,297)
    ,(376,6,15)
    ,(376,7,12)
    ,(376,8,9)
    ,(376,9,-4)
    ,(376,10,-7)
    ,(376,11,-12)
    ,(376,12,35)
    ,(376,13,25)
    ,(376,14,7)
    ,(376,15,1015)
    ,(376,16,1013)
    ,(376,17,1006)
    ,(376,18,14)
    ,(376,19,11)
    ,(376,20,10)
    ,(376,21,55)
    ,(376,22,29)
    ,(376,24,0)
    ,(376,25,1)
    ,(376,27,280)
    ,(377,6,17)
    ,(377,7,11)
    ,(377,8,6)
    ,(377,9,-1)
    ,(377,10,-7)
    ,(377,11,-12)
    ,(377,12,46)
    ,(377,13,28)
    ,(377,14,6)
    ,(377,15,1018)
    ,(377,16,1014)
    ,(377,17,1010)
    ,(377,18,11)
    ,(377,19,10)
    ,(377,20,10)
    ,(377,21,26)
    ,(377,22,10)
    ,(377,24,0)
    ,(377,25,4)
    ,(377,27,267)
    ,(378,6,7)
    ,(378,7,6)
    ,(378,8,5)
    ,(378,9,4)
    ,(378,10,2)
    ,(378,11,-3)
    ,(378,12,87)
    ,(378,13,71)
    ,(378,14,33)
    ,(378,15,1021)
    ,(378,16,1016)
    ,(378,17,1011)
    ,(378,18,10)
    ,(378,19,8)
    ,(378,20,3)
    ,(378,21,26)
    ,(378,22,16)
    ,(378,24,4.0599999999999996)
    ,(378,25,7)
    ,(378,26,"Rain")
    ,(378,27,238)
    ,(379,6,14)
    ,(379,7,8)
    ,(379,8,2)
    ,(379,9,1)
    ,(379,10,-3)
    ,(379,11,-6)
    ,(379,12,81)
    ,(379,13,47)
    ,(379,14,14)
    ,(379,15,1021)
    ,(379,16,1018)
    ,(379,17,1016)
    ,(379,18,14)
    ,(379,19,10)
    ,(379,20,10)
    ,(379,21,34)
    ,(379,22,14)
    ,(379,24,0)
    ,(379,25,2)
    ,(379,27,268)
    ,(380,6,18)
    ,(380,7,12)
    ,(380,8,6)
    ,(380,9,-2)
    ,(380,10,-7)
    ,(380,11,-15)
    ,(380,12,57)
    ,(380,13,28)
    ,(380,14,5)
    ,(380,15,1018)
    ,(380,16,1016)
    ,(380,17,1012)
    ,(380,18,14)
    ,(380,19,10)
    ,(380,20,10)
    ,(380,21,19)
    ,(380,22,6)
    ,(380,24,0)
    ,(380,25,4)
    ,(380,27,268)
    ,(381,6,19)
    ,(381,7,14)
    ,(381,8,11)
    ,(381,9,2)
    ,(381,10,-3)
    ,(381,11,-7)
    ,(381,12,54)
    ,(381,13,29)
    ,(381,14,12)
    ,(381,15,1015)
    ,(381,16,1014)
    ,(381,17,1008)
    ,(381,18,14)
    ,(381,19,10)
    ,(381,20,10)
    ,(381,21,48)
    ,(381,22,18)
    ,(381,23,55)
    ,(381,24,0)
    ,(381,25,3)
    ,(381,27,280)
    ,(382,6,21)
    ,(382,7,15)
    ,(382,8,9)
    ,(382,9,7)
    ,(382,10,-3)
    ,(382,11,-9)
    ,(382,12,68)
    ,(382,13,30)
    ,(382,14,8)
    ,(382,15,1015)
    ,(382,16,1013)
    ,(382,17,1008)
    ,(382,18,14)
    ,(382,19,10)
    ,(382,20,10)
    ,(382,21,37)
    ,(382,22,19)
    ,(382,24,0)
    ,(382,25,2)
    ,(382,27,274)
    ,(383,6,22)
    ,(383,7,17)
    ,(383,8,11)
    ,(383,9,0)
    ,(383,10,-4)
    ,(383,11,-8)
    ,(383,12,44)
    ,(383,13,23)
    ,(383,14,7)
    ,(383,15,1016)
    ,(383,16,1013)
    ,(383,17,1008)
    ,(383,18,11)
    ,(383,19,10)
    ,(383,20,8)
    ,(383,21,19)
    ,(383,22,10)
    ,(383,24,0)
    ,(383,25,1)
    ,(383,27,247)
    ,(384,6,24)
    ,(384,7,18)
    ,(384,8,12)
    ,(384,9,-3)
    ,(384,10,-6)
    ,(384,11,-11)
    ,(384,12,33)
    ,(384,13,18)
    ,(384,14,4)
    ,(384,15,1016)
    ,(384,16,1013)
    ,(384,17,1007)
    ,(384,18,11)
    ,(384,19,10)
    ,(384,20,10)
    ,(384,21,37)
    ,(384,22,11)
    ,(384,24,0)
    ,(384,25,1)
    ,(384,27,290)
    ,(385,6,23)
    ,(385,7,17)
    ,(385,8,10)
    ,(385,9,3)
    ,(385,10,-2)
    ,(385,11,-5)
    ,(385,12,45)
    ,(385,13,27)
    ,(385,14,8)
    ,(385,15,1018)
    ,(385,16,1015)
    ,(385,17,1009)
    ,(385,18,11)
    ,(385,19,10)
    ,(385,20,10)
    ,(385,21,19)
    ,(385,22,8)
    ,(385,23,37)
    ,(385,24,0)
    ,(385,25,3)
    ,(385,27,288)
    ,(386,6,20)
    ,(386,7,16)
    ,(386,8,11)
    ,(386,9,7)
    ,(386,10,3)
    ,(386,11,0)
    ,(386,12,72)
    ,(386,13,40)
    ,(386,14,21)
    ,(386,15,1019)
    ,(386,16,1016)
    ,(386,17,1011)
    ,(386,18,10)
    ,(386,19,10)
    ,(386,20,7)
    ,(386,21,23)
    ,(386,22,10)
    ,(386,24,1.02)
    ,(386,25,4)
    ,(386,26,"Rain")
    ,(386,27,90)
    ,(387,6,23)
    ,(387,7,18)
    ,(387,8,12)
    ,(387,9,10)
    ,(387,10,3)
    ,(387,11,-1)
    ,(387,12,68)
    ,(387,13,37)
    ,(387,14,15)
    ,(387,15,1019)
    ,(387,16,1017)
    ,(387,17,1011)
    ,(387,18,10)
    ,(387,19,10)
    ,(387,20,10)
    ,(387,21,34)
    ,(387,22,8)
    ,(387,24,0)
    ,(387,25,4)
    ,(387,26,"Rain")
    ,(387,27,261)
    ,(388,6,23)
    ,(388,7,17)
    ,(388,8,11)
    ,(388,9,10)
    ,(388,10,7)
    ,(388,11,2)
    ,(388,12,82)
    ,(388,13,56)
    ,(388,14,17)
    ,(388,15,1017)
    ,(388,16,1015)
    ,(388,17,1009)
    ,(388,18,10)
    ,(388,19,8)
    ,(388,20,3)
    ,(388,21,40)
    ,(388,22,16)
    ,(388,24,2.0299999999999998)
    ,(388,25,6)
    ,(388,26,"Rain")
    ,(388,27,268)
)

foreach ($entry in $cellData) {
    $r = $entry[0]
    $c = $entry[1]
    $v = $entry[2]
    $ws.Cells.Item($r, $c).Value = $v
}

$ws.Range("C377").Select()

Write-Host "Added rows 366-388 to Sheet1"